$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure static columns A,B,C,E,F,G,H,I,N,O,Q,R are set for row 25 (new row),
# copying the constant values used throughout the sheet.
$ws.Range("A25").Value2 = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("E25").Value2 = 9
$ws.Range("F25").Value2 = 100112042
$ws.Range("G25").Value = "Locoto"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("N25").Value = "$/kilo"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("Q25").Value2 = 1
$ws.Range("R25").Value = "Hortaliza"

# The date column uses a custom date number format throughout the sheet;
# apply the same format to the newly added row so D25 renders as a date.
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update the weekly price/volume/date data for rows 2-25 (values shifted down
# by one row, with a new most-recent entry placed at row 2).
$ws.Range("D2").Value2 = 44749
$ws.Range("J2").Value2 = 80
$ws.Range("K2").Value2 = 2500
$ws.Range("L2").Value2 = 2500
$ws.Range("M2").Value2 = 2500
$ws.Range("P2").Value2 = 2500

$ws.Range("D3").Value2 = 44677
$ws.Range("J3").Value2 = 20
$ws.Range("K3").Value2 = 5500
$ws.Range("L3").Value2 = 5500
$ws.Range("M3").Value2 = 5500
$ws.Range("P3").Value2 = 5500

$ws.Range("D4").Value2 = 44719
$ws.Range("J4").Value2 = 80
$ws.Range("K4").Value2 = 3600
$ws.Range("L4").Value2 = 3600
$ws.Range("M4").Value2 = 3600
$ws.Range("P4").Value2 = 3600

$ws.Range("D5").Value2 = 44203
$ws.Range("J5").Value2 = 30
$ws.Range("K5").Value2 = 2000
$ws.Range("L5").Value2 = 2000
$ws.Range("M5").Value2 = 2000
$ws.Range("P5").Value2 = 2000

$ws.Range("D6").Value2 = 44669
$ws.Range("J6").Value2 = 60
$ws.Range("K6").Value2 = 6250
$ws.Range("L6").Value2 = 6250
$ws.Range("M6").Value2 = 6250
$ws.Range("P6").Value2 = 6250

$ws.Range("D7").Value2 = 44497
$ws.Range("J7").Value2 = 50
$ws.Range("K7").Value2 = 2200
$ws.Range("L7").Value2 = 2200
$ws.Range("M7").Value2 = 2200
$ws.Range("P7").Value2 = 2200

$ws.Range("D8").Value2 = 44740
$ws.Range("J8").Value2 = 50
$ws.Range("K8").Value2 = 2500
$ws.Range("L8").Value2 = 2500
$ws.Range("M8").Value2 = 2500
$ws.Range("P8").Value2 = 2500

$ws.Range("D9").Value2 = 44741
$ws.Range("J9").Value2 = 100
$ws.Range("K9").Value2 = 2500
$ws.Range("L9").Value2 = 2500
$ws.Range("M9").Value2 = 2500
$ws.Range("P9").Value2 = 2500

$ws.Range("D10").Value2 = 44679
$ws.Range("J10").Value2 = 30
$ws.Range("K10").Value2 = 5500
$ws.Range("L10").Value2 = 5500
$ws.Range("M10").Value2 = 5500
$ws.Range("P10").Value2 = 5500

$ws.Range("D11").Value2 = 44447
$ws.Range("J11").Value2 = 75
$ws.Range("K11").Value2 = 2200
$ws.Range("L11").Value2 = 2200
$ws.Range("M11").Value2 = 2200
$ws.Range("P11").Value2 = 2200

$ws.Range("D12").Value2 = 44453
$ws.Range("J12").Value2 = 20
$ws.Range("K12").Value2 = 2300
$ws.Range("L12").Value2 = 2300
$ws.Range("M12").Value2 = 2300
$ws.Range("P12").Value2 = 2300

$ws.Range("D13").Value2 = 44487
$ws.Range("J13").Value2 = 50
$ws.Range("K13").Value2 = 2200
$ws.Range("L13").Value2 = 2200
$ws.Range("M13").Value2 = 2200
$ws.Range("P13").Value2 = 2200

$ws.Range("D14").Value2 = 44496
$ws.Range("J14").Value2 = 40
$ws.Range("K14").Value2 = 2200
$ws.Range("L14").Value2 = 2200
$ws.Range("M14").Value2 = 2200
$ws.Range("P14").Value2 = 2200

$ws.Range("D15").Value2 = 44484
$ws.Range("J15").Value2 = 40
$ws.Range("K15").Value2 = 2200
$ws.Range("L15").Value2 = 2200
$ws.Range("M15").Value2 = 2200
$ws.Range("P15").Value2 = 2200

$ws.Range("D16").Value2 = 44685
$ws.Range("J16").Value2 = 60
$ws.Range("K16").Value2 = 5000
$ws.Range("L16").Value2 = 6000
$ws.Range("M16").Value2 = 5333
$ws.Range("P16").Value2 = 5333

$ws.Range("D17").Value2 = 44720
$ws.Range("J17").Value2 = 100
$ws.Range("K17").Value2 = 3600
$ws.Range("L17").Value2 = 3600
$ws.Range("M17").Value2 = 3600
$ws.Range("P17").Value2 = 3600

$ws.Range("D18").Value2 = 44452
$ws.Range("J18").Value2 = 120
$ws.Range("K18").Value2 = 2300
$ws.Range("L18").Value2 = 2300
$ws.Range("M18").Value2 = 2300
$ws.Range("P18").Value2 = 2300

$ws.Range("D19").Value2 = 44706
$ws.Range("J19").Value2 = 90
$ws.Range("K19").Value2 = 4700
$ws.Range("L19").Value2 = 4700
$ws.Range("M19").Value2 = 4700
$ws.Range("P19").Value2 = 4700

$ws.Range("D20").Value2 = 44476
$ws.Range("J20").Value2 = 30
$ws.Range("K20").Value2 = 2200
$ws.Range("L20").Value2 = 2200
$ws.Range("M20").Value2 = 2200
$ws.Range("P20").Value2 = 2200

$ws.Range("D21").Value2 = 44474
$ws.Range("J21").Value2 = 20
$ws.Range("K21").Value2 = 1600
$ws.Range("L21").Value2 = 1600
$ws.Range("M21").Value2 = 1600
$ws.Range("P21").Value2 = 1600

$ws.Range("D22").Value2 = 44473
$ws.Range("J22").Value2 = 140
$ws.Range("K22").Value2 = 1600
$ws.Range("L22").Value2 = 1600
$ws.Range("M22").Value2 = 1600
$ws.Range("P22").Value2 = 1600

$ws.Range("D23").Value2 = 44707
$ws.Range("J23").Value2 = 100
$ws.Range("K23").Value2 = 4700
$ws.Range("L23").Value2 = 4700
$ws.Range("M23").Value2 = 4700
$ws.Range("P23").Value2 = 4700

$ws.Range("D24").Value2 = 44747
$ws.Range("J24").Value2 = 80
$ws.Range("K24").Value2 = 2500
$ws.Range("L24").Value2 = 2500
$ws.Range("M24").Value2 = 2500
$ws.Range("P24").Value2 = 2500

$ws.Range("D25").Value2 = 44483
$ws.Range("J25").Value2 = 50
$ws.Range("K25").Value2 = 2200
$ws.Range("L25").Value2 = 2200
$ws.Range("M25").Value2 = 2200
$ws.Range("P25").Value2 = 2200

